# multi browser implementation for chrome and firefox and accounting cash
#
# Update the "currency" row on the ProductLoanInput sheet (label + value),
# update the short-name value on the ProductLoanOutput sheet, and make
# ProductLoanInput the active sheet/selection (was ProductLoanOutput).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# --- ProductLoanInput: row 6 (currency / US Dollar) -----------------------
$ws1.Range("A6").Value = "currency"

$ws1.Range("B6").Value = "US Dollar"
# B6 picks up the same look as the other plain "green fill" input cells
# (style index 1 in the original workbook) instead of the wrapped/indented
# style it had before.
$ws1.Range("B6").Interior.Color = 5296274
$ws1.Range("B6").Font.Name = "Arial"
$ws1.Range("B6").Font.Size = 10

# --- ProductLoanOutput: B1 short name --------------------------------------
$ws2.Range("B1").Value = "823RBI-EI-DB-SAR-Y-NON-RNI-CTRFD-SAR-MD-Y-1-Late Repayment"

# --- Active sheet / selection ----------------------------------------------
# Previously ProductLoanOutput was the active tab with B1 selected and
# ProductLoanInput was scrolled to D7; now ProductLoanInput is the active
# tab with A6:B6 selected (scroll position reset).
$ws1.Activate()
$ws1.Range("A6:B6").Select()

Write-Output "edit complete"
